# INCLUIDO COD DO EVENTO
# Replace the VLOOKUP-against-external-workbook formulas in column C
# (rows 2-26) with their resolved static values, then break the now-unused
# external link so the external link parts are removed from the package.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Set every cell C2:C26 to 27, the value shared by almost all rows...
$ws.Range("C2:C26").Value = 27

# ...except C18, which resolves to 20.
$ws.Range("C18").Value = 20

# The formulas were the only thing referencing the external workbook link;
# breaking it removes xl/externalLinks/externalLink1.xml, the
# <externalReference> entry in workbook.xml, and the now-empty calcChain.
foreach ($source in $wb.LinkSources()) {
    $wb.BreakLink($source, 1)
}

# Move the active selection to E18, matching the author's final cursor spot.
$ws.Range("E18").Select()
